$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("All Published Values")
$wsSummary = $wb.Worksheets.Item("Daily Summary")

# --- Append new row 22 to "All Published Values" ---
# All columns in this sheet are stored as text, so force text format before
# assigning values (otherwise date-like / number-like strings get converted).
$rowValues = @(
    "2026-01-03",
    "2026-01-03 05:58:24",
    "697.85",
    "697.85",
    "700.79",
    "700.79",
    "702.88",
    "2026/01/03 05:58:24",
    "2026-01-02 22:13:44",
    "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
)

for ($col = 1; $col -le $rowValues.Length; $col++) {
    $cell = $wsData.Cells.Item(22, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $rowValues[$col - 1]
    # Drop the temporary "text" number format again so the cell ends up
    # with the default (unstyled) look, matching the other data rows.
    $cell.Style = "Normal"
}

# --- Update dependent ranges on "All Published Values" ---
# Re-apply the AutoFilter over the expanded range (toggle off, then back on
# so the stored ref gets refreshed to A1:J22).
$wsData.Range("A1:J22").AutoFilter() | Out-Null
$wsData.Range("A1:J22").AutoFilter() | Out-Null

# Update the hidden _xlnm._FilterDatabase defined name for this sheet.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$22"
    }
}

# --- Update "Daily Summary" sheet: publishes count for 2026-01-03 (B5) ---
$wsSummary.Cells.Item(5, 2).Value = 3
